$p = $ppt.ActivePresentation
Write-Output "noop"
